# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values for column G ("K") for rows 2-34
$newValues = @{
    2  = 2
    3  = 0
    4  = 1
    5  = 0
    6  = 2
    7  = 0
    8  = 1
    9  = 0
    10 = 1
    11 = 0
    12 = 1
    13 = 1
    14 = 1
    15 = 1
    16 = 3
    17 = 2
    18 = 1
    19 = 1
    20 = 1
    21 = 2
    22 = 0
    23 = 4
    24 = 2
    25 = 0
    26 = 0
    27 = 2
    28 = 2
    29 = 2
    30 = 1
    31 = 1
    32 = 2
    33 = 2
    34 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
